$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Settings"

# Update row 2 (Manufacturer_Path) value
$ws.Range("B2").Value = "Data\MasterULRs.xlsx"

# Update row 3: was Start_URL / https://www.google.com/, becomes ReportPath / D:\ReportsScraping\
$ws.Range("A3").Value = "ReportPath"
$ws.Range("B3").Value = "D:\ReportsScraping\"

# Remove old row 4 (ReportPath / D:\ReportsScraping\) which is now redundant
$ws.Range("A4:B4").Delete()

# Update selection to B8 as per diff
$ws.Range("B8").Select()

# Update window view size/position (maximized window geometry)
$win = $wb.Windows.Item(1)
$win.Left = -108
$win.Top = -108
$win.Width = 23256
$win.Height = 12456
$win.WindowState = -4137
